# rename flow-of-control to control flow
$p = $ppt.ActivePresentation

# Slide 29: "The Flow-Of-Control is altered ..." -> "The control flow is altered ..."
$s29 = $p.Slides.Item(29)
$sh29 = $s29.Shapes.Item(2)
$para29 = $sh29.TextFrame.TextRange.Paragraphs(1, 1)
# Use a two-step write so the run is fully replaced (not diff-merged into
# several runs) and collapses to a single run like the target XML.
$para29.Text = "PLACEHOLDER"
$para29b = $sh29.TextFrame.TextRange.Paragraphs(1, 1)
$para29b.Text = "The control flow is altered " + [char]0x2013 + " methods can return when an await statement is executed"

# Slide 5: merge the three "12KB (" / "x86) or 14 " / "KB (x64) Kernel Mode Object" runs
# into a single run with unchanged visible text.
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$para5 = $sh5.TextFrame.TextRange.Paragraphs(4, 1)
$para5.Text = "PLACEHOLDER"
$para5b = $sh5.TextFrame.TextRange.Paragraphs(4, 1)
$para5b.Text = "12KB (x86) or 14 KB (x64) Kernel Mode Object"
